$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestSteps")

# A new "checkAccessibility" test step needs to become row 3, pushing the
# existing rows 3..12 down to 4..13. Rows.Insert() in this runtime always
# stamps the freshly-created row with a generic "blank row" style (and
# mints a brand-new, otherwise-unused cellXf for it), which diverges from
# the target file (every data row keeps using the original style index 3,
# exactly as styles.xml already defines). So instead the rows are shifted
# manually, bottom-up, by copying formats then values into the row below.
#
# Only the original row 3 carries values in columns D:F (a quirk inherited
# from the sheet's very first row), every other data row only uses A:C -
# match that per row so no stray empty D/E/F cells get created.
for ($r = 12; $r -ge 3; $r--) {
    if ($r -eq 3) {
        $lastCol = "F"
    } else {
        $lastCol = "C"
    }
    $src = $ws.Range("A" + $r + ":" + $lastCol + $r)
    $dst = $ws.Range("A" + ($r + 1) + ":" + $lastCol + ($r + 1))
    $src.Copy()
    $dst.PasteSpecial(-4122)
    $src.Copy()
    $dst.PasteSpecial(-4163)
}

$ws.Cells.Item(3, 1).Value = "checkAccessibility"
$ws.Cells.Item(3, 2).Value = "AddNewProfile_OrganizationCodeProfile"
$ws.Cells.Item(3, 3).ClearContents()

# Column B now holds the longer ObjectID text above, so it widens to fit.
$ws.Columns("B").ColumnWidth = 34.46

# Leave the same kind of single-cell selection state the workbook was
# saved with, just shifted one column left.
$null = $ws.Range("C12").Select()
